$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Give the "Encrypted_files" sub-table a new trailing column (D).
#     Seed it (values + formatting) from the sheet's old last column (K,
#     "auth_image"), which already carries the correct "final column"
#     border treatment (header box + left/right body borders). Do this
#     first, while K still holds its original content. ---
$ws.Range("K5:K8").Copy()
$ws.Range("D5").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("K5:K8").Copy()
$ws.Range("D5").PasteSpecial(-4163)   # xlPasteValues
$ws.Application.CutCopyMode = $false

# --- Drop the stray duplicate "pid" column (E) that used to sit between
#     the two sub-tables; the gap column stays, just emptied out. ---
$ws.Range("E5:E8").Clear()

# --- Shift the "User" sub-table left by one column, dropping the old
#     "email" column: username stays in F, and four_digit_pass /
#     private_key / profile_pic each slide one column left (G<-H,
#     H<-I, I<-J), copying values+formats left-to-right so a donor
#     column is always read before it gets overwritten. ---
$ws.Range("H5:H8").Copy()
$ws.Range("G5").PasteSpecial(-4122)
$ws.Range("H5:H8").Copy()
$ws.Range("G5").PasteSpecial(-4163)

$ws.Range("I5:I8").Copy()
$ws.Range("H5").PasteSpecial(-4122)
$ws.Range("I5:I8").Copy()
$ws.Range("H5").PasteSpecial(-4163)

$ws.Range("J5:J8").Copy()
$ws.Range("I5").PasteSpecial(-4122)
$ws.Range("J5:J8").Copy()
$ws.Range("I5").PasteSpecial(-4163)

$ws.Application.CutCopyMode = $false

# --- Clear the now-vacated trailing columns (old profile_pic / auth_image
#     slots) so the "User" sub-table ends cleanly at I. ---
$ws.Range("J5:J8").Clear()
$ws.Range("K5:K8").Clear()

# --- Rename the headers that changed wording. ---
$ws.Range("C5").Value = "Filename"
$ws.Range("I5").Value = "profile_pic_filename"
$ws.Range("A5").Value = "id"
$ws.Range("D5").Value = "date_uploaded"

# --- Match the saved selection. ---
$ws.Range("D12").Select()
